$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change %) per upstream diff.
# D-column prices are plain text (thousand-dot formatted, e.g. '70.012.90'),
# so force Text number format before assignment to avoid Excel's automatic
# numeric coercion (which would drop meaningful trailing zeros / thousand dots).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.973.61'
$ws.Range("E2").Value = '  -4.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.806.32'
$ws.Range("E3").Value = '  -4.87%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.47'
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.94'
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.660'
$ws.Range("E7").Value = '  -3.78%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.08'
$ws.Range("E11").Value = '  -4.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000316'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.406.50'
$ws.Range("E14").Value = '  -4.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.796.00'
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.65'
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("E17").Value = '  -4.21%  '
$ws.Range("E18").Value = '  -6.57%  '
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.789.17'
$ws.Range("E20").Value = '  -4.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.49'
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.65'
$ws.Range("E22").Value = '  -3.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '92.50'
$ws.Range("E23").Value = '  -4.18%  '
$ws.Range("E24").Value = '  -6.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.67'
$ws.Range("E25").Value = '  -4.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.05'
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.91'
$ws.Range("E27").Value = '  -12.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.95'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.31'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.64'
$ws.Range("E30").Value = '  -4.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.01'
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.30'
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '47.66'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("E34").Value = '  -6.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '68.52'
$ws.Range("E35").Value = '  -3.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0962'
$ws.Range("E36").Value = '  +9.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '625.26'
$ws.Range("E37").Value = '  -7.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.418'
$ws.Range("E38").Value = '  -5.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("E42").Value = '  -4.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.17'
$ws.Range("E43").Value = '  +19.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0462'
$ws.Range("E44").Value = '  -5.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.73'
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.79'
$ws.Range("E46").Value = '  -9.32%  '
$ws.Range("E47").Value = '  -5.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").Value = '  -16.09%  '
$ws.Range("E49").Value = '  -5.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.793.32'
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("E51").Value = '  -0.38%  '
